$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header changes
$ws.Range("G1").Value = "S Tag"

# D3: Student -> S, and related G3 text
$ws.Range("D3").Value = "S"
$ws.Range("G3").Value = "2 - relating to another S"

# D6: Students -> SS
$ws.Range("D6").Value = "SS"

# D column: RBD -> T for all matching rows
$rbdRows = @(2,7,9,12,14,15,17,25,26,28,29,30,34,47,48,50,52,54,56,59,61,62,63,65,66,68,69,72,74,75,80,81,82,83,84,85,86,87,88)
foreach ($r in $rbdRows) {
    $ws.Range("D$r").Value = "T"
}
